$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely (removes "Matrix Generation Time" row)
$ws.Rows.Item(2).Delete()

# Update remaining row's values
$ws.Range("A1").Value = "Parallel Multiplication Time"
$ws.Range("B1").Value = 0.0300042
